$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 29   Number  47"
$ws.Range("C9").Value = "Report Covering the Week  11/21/2022  Through  11/27/2022"

# --- Simple numeric value updates (style unchanged) ---
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = 66.666666666666
$ws.Range("I16").Value = 49
$ws.Range("K16").Value = 4.255319148936
$ws.Range("L16").Value = -25.757575757575
$ws.Range("M16").Value = -48.421052631578
$ws.Range("N16").Value = -90.576923076923
$ws.Range("F17").Value = 4
$ws.Range("G17").Value = 6
$ws.Range("H17").Value = -33.333333333333
$ws.Range("I17").Value = 92
$ws.Range("K17").Value = 13.580246913580
$ws.Range("L17").Value = -26.4
$ws.Range("M17").Value = -3.157894736842
$ws.Range("N17").Value = -64.885496183206
$ws.Range("C18").Value = 3
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = -61.538461538461
$ws.Range("I18").Value = 88
$ws.Range("K18").Value = -23.478260869565
$ws.Range("L18").Value = -30.158730158730
$ws.Range("M18").Value = -61.233480176211
$ws.Range("N18").Value = -91.423001949317
$ws.Range("C19").Value = 10
$ws.Range("E19").Value = 42.857142857142
$ws.Range("F19").Value = 37
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = 32.142857142857
$ws.Range("I19").Value = 460
$ws.Range("J19").Value = 335
$ws.Range("K19").Value = 37.313432835820
$ws.Range("L19").Value = 55.405405405405
$ws.Range("M19").Value = 55.405405405405
$ws.Range("N19").Value = 8.490566037735
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = -40
$ws.Range("I20").Value = 96
$ws.Range("J20").Value = 81
$ws.Range("K20").Value = 18.518518518518
$ws.Range("L20").Value = -21.951219512195
$ws.Range("M20").Value = -30.434782608695
$ws.Range("N20").Value = -94.606741573033
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 12
$ws.Range("E21").Value = 50
$ws.Range("F21").Value = 57
$ws.Range("G21").Value = 61
$ws.Range("H21").Value = -6.557377049180
$ws.Range("I21").Value = 797
$ws.Range("J21").Value = 672
$ws.Range("K21").Value = 18.601190476190
$ws.Range("L21").Value = 5.562913907284
$ws.Range("M21").Value = -7.109557109557
$ws.Range("N21").Value = -80.291790306627
$ws.Range("D24").Value = 38
$ws.Range("E24").Value = 10.526315789473
$ws.Range("F24").Value = 149
$ws.Range("G24").Value = 128
$ws.Range("H24").Value = 16.40625
$ws.Range("I24").Value = 1715
$ws.Range("J24").Value = 996
$ws.Range("K24").Value = 72.188755020080
$ws.Range("L24").Value = 90.344062153163
$ws.Range("M24").Value = 56.050955414012
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 25
$ws.Range("F25").Value = 22
$ws.Range("G25").Value = 16
$ws.Range("H25").Value = 37.5
$ws.Range("I25").Value = 288
$ws.Range("J25").Value = 224
$ws.Range("K25").Value = 28.571428571428
$ws.Range("L25").Value = 35.849056603773
$ws.Range("M25").Value = -19.777158774373
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = -40
$ws.Range("J27").Value = 35
$ws.Range("K27").Value = 20
$ws.Range("L27").Value = 10.526315789473

# --- Text -> Number conversions (copy number-style format, then set value) ---
$ws.Range("I15").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("D14").Value = 1
$ws.Range("M15").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("E14").Value = -100
$ws.Range("I15").Copy()
$ws.Range("G14").PasteSpecial(-4122)
$ws.Range("G14").Value = 1
$ws.Range("M15").Copy()
$ws.Range("H14").PasteSpecial(-4122)
$ws.Range("H14").Value = -100
$ws.Range("I15").Copy()
$ws.Range("J14").PasteSpecial(-4122)
$ws.Range("J14").Value = 1
$ws.Range("M15").Copy()
$ws.Range("K14").PasteSpecial(-4122)
$ws.Range("K14").Value = 100

# --- Number -> Text conversions (force text entry via quote-prefix, then copy text-style format) ---
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0"
$ws.Range("C23").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "***.*"
$ws.Range("C23").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0"
$ws.Range("C23").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "***.*"
$ws.Range("C23").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0"
$ws.Range("C23").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "***.*"
$ws.Range("C23").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("C23").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("C23").Copy()
$ws.Range("C27").PasteSpecial(-4122)

Write-Host "Edit complete"